# Weekly update: insert a new data row at row 9 (pushing existing rows 9-13
# down to 10-14) with the latest week's price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9, shifting rows 9:13 down to 10:14.
$ws.Rows(9).Insert()

# Populate the newly inserted row 9 with the new observation. The row shares
# the same constant columns (mercado/region/category/etc.) as the rest of the
# dataset, only the date/volume/price/origin/unit columns change.
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44704
$ws.Cells.Item(9, 4).Style = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(10, 4).NumberFormat
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = 100112037
$ws.Cells.Item(9, 7).Value = "Cebollín"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 6000
$ws.Cells.Item(9, 12).Value = 6500
$ws.Cells.Item(9, 13).Value = 6250
$ws.Cells.Item(9, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(9, 15).Value = "Región Metropolitana"
$ws.Cells.Item(9, 16).Value = 174
$ws.Cells.Item(9, 17).Value = 36
$ws.Cells.Item(9, 18).Value = "Hortaliza"
